$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timecard hours:
#  Andrew Case (Mine)   +4 hours : 29h 29m -> 33h 29m
#  Matthew Darby (Matt) +2 hours : 16h 15m -> 18h 15m
#  Weston Straw         +4 hours : 19h 14m -> 23h 14m
$ws.Range("B4").Value = "33h 29m"
$ws.Range("B9").Value = "23h 14m"
$ws.Range("B5").Value = "18h 15m"

# Move the active selection to B5
$ws.Range("B5").Select()
